{"js": "// Apply hybrid bold + color highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific resume\n// bullet paragraphs, matching the target diff exactly.\n//\n// Strategy: for each target paragraph (identified by its exact, unique\n// full paragraph text), run an in-paragraph search() for each metric\n// substring (in left-to-right order) and apply bold + the brand color\n// (#2C3E50) to just that sub-range. Word/Office.js's search() naturally\n// splits the existing run(s) around the match, which reproduces the\n// run-split pattern shown in the diff (plain-text run, bold+colored\n// metric run, plain-text run, ...).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Exact full-text -> ordered list of metric substrings to highlight.\nconst EDITS = [\n  {\n    text: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    text: \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nfor (const edit of EDITS) {\n  // Locate the paragraph with this exact text (each target string is\n  // unique in the document, so the first/only match is the right one).\n  const para = paragraphs.items.find((p) => p.text === edit.text);\n  if (!para) continue;\n\n  for (const metric of edit.metrics) {\n    const found = para.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < found.items.length; i++) {\n      const rng = found.items[i];\n      rng.font.bold = true;\n      rng.font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific resume\n# bullet paragraphs, matching the target diff exactly.\n#\n# Strategy: for each target paragraph (identified by its exact, unique\n# full paragraph text, trailing paragraph-mark trimmed), duplicate the\n# paragraph's Range and run Find.Execute for each metric substring (in\n# left-to-right order), then apply Bold + the brand color (#2C3E50) to\n# just that found sub-range. Word's Find/Range model naturally splits\n# the existing run(s) around the match, reproducing the run-split\n# pattern shown in the diff (plain-text run, bold+colored metric run,\n# plain-text run, ...).\n\nfunction Get-WdColor([int]$r, [int]$g, [int]$b) {\n    # Word COM Font.Color is a 24-bit BGR-packed integer (R + G*256 + B*65536),\n    # not the RRGGBB order used by OOXML's <w:color w:val=\"RRGGBB\"/>.\n    return $r + ($g * 256) + ($b * 65536)\n}\n\n$HighlightColor = Get-WdColor 0x2C 0x3E 0x50\n\n$d = $word.ActiveDocument\n\n$Edits = @(\n    @{\n        Text    = ([char]0x2022 + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\")\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = ([char]0x2022 + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + [char]0xB1 + \"4.2% to \" + [char]0xB1 + \"2.1%\")\n        Metrics = @(\"87%\", \"71%\", ([char]0xB1 + \"4.2%\"), ([char]0xB1 + \"2.1%\"))\n    },\n    @{\n        Text    = ([char]0x2022 + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\")\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Text    = ([char]0x2022 + \" Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\")\n        Metrics = @('$400M', '$1B')\n    },\n    @{\n        Text    = ([char]0x2022 + \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\")\n        Metrics = @(\"73.5%\", '$4.7M')\n    },\n    @{\n        Text    = ([char]0x2022 + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\")\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    # Word COM's Range.Text includes the trailing paragraph mark (\\r) -\n    # trim it so we can compare against the plain diff text.\n    if ($t.Length -gt 0 -and [int][char]$t[$t.Length - 1] -eq 13) {\n        $t = $t.Substring(0, $t.Length - 1)\n    }\n\n    foreach ($edit in $Edits) {\n        if ($t -eq $edit.Text) {\n            $pRange = $p.Range\n            foreach ($metric in $edit.Metrics) {\n                $rng = $pRange.Duplicate\n                $rng.Find.Text = $metric\n                $found = $rng.Find.Execute()\n                if ($found) {\n                    $rng.Font.Bold = 1\n                    $rng.Font.Color = $HighlightColor\n                }\n            }\n        }\n    }\n}\n"}
